# Adds a leading "Posición" (row rank) column and a trailing "CPK" column
# to the two "Top 10 Rutas" summary sheets.

$wb = $excel.ActiveWorkbook

$cpkBySheet = @{
    "Top 10 Rutas Mas Eficientes" = @(19.31470044841836, 6.470371847307608, 22.94552710008876, 6.043426602409495, 23.78267234590579, 27.8651583241622, 11.88791422506719, 27.00580232297981, 7.779807821051387, 7.283936170943039)
    "Top 10 Rutas Menos Eficientes" = @(26742.72581176471, 19118.89281191969, 1107.15137569079, 4570.410677494797, 5481.384329607672, 6157.743728477009, 6039.604892662821, 2176.151316410839, 2645.576421704518, 4989.200344019353)
}

foreach ($sheetName in @("Top 10 Rutas Mas Eficientes", "Top 10 Rutas Menos Eficientes")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert new column A ("Posición"), shifting Ruta..Eficiencia from A-E to B-F.
    $ws.Columns.Item(1).Insert()

    # Copy the header style (bold/centered/bordered) from the old-A (now B) header onto the new A1.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Posición"

    for ($i = 2; $i -le 11; $i++) {
        $ws.Cells.Item($i, 1).Value = $i - 1
    }

    # Append new column G ("CPK"), copying header style from column F.
    $ws.Range("F1").Copy()
    $ws.Range("G1").PasteSpecial(-4122)
    $ws.Range("G1").Value = "CPK"

    $cpk = $cpkBySheet[$sheetName]
    for ($i = 0; $i -le 9; $i++) {
        $ws.Cells.Item($i + 2, 7).Value = $cpk[$i]
    }
}
